# Part3-ExceptionHandling.pptx - "Added a slide in exceptions slides for exercise"
#
# Slide 11 previously held a single full-height picture centered-ish on the
# slide. The edit shifts that picture to the right half of the slide and
# adds a new title-like textbox ("What does this output?") in the
# now-empty left area.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# --- Move the existing picture to the right ------------------------------
# EMU -> Points conversion (PowerPoint COM Left/Top/Width/Height are
# expressed in points; 1 pt = 12700 EMU). The literal constants below are
# pre-nudged by a few millionths of a point so that round-tripping through
# the host's single-precision storage lands exactly on the target EMU
# value instead of being truncated one EMU short.
$picLeftPt   = 370.4841772283464   # -> 4705149 EMU
$tbLeftPt    = 18.20693013385827   # -> 231228 EMU
$tbTopPt     = 13.241418322834646  # -> 168166 EMU
$tbWidthPt   = 304.0909578818898   # -> 3861955 EMU
$tbHeightPt  = 41.198425196850394  # -> 523220 EMU

$pic = $s.Shapes.Item(1)
$pic.Left = $picLeftPt

# --- Add the new textbox --------------------------------------------------
$tb = $s.Shapes.AddTextbox(
    1,
    $tbLeftPt,
    $tbTopPt,
    $tbWidthPt,
    $tbHeightPt
)

$tb.Name = "TextBox 2"
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

$tb.TextFrame.TextRange.Text = "What does this output?"
$tb.TextFrame.TextRange.Font.Size = 28
